$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill A1:A6 with the value 1
$ws.Range("A1:A6").Value = 1

# Leave the selection where the user ended up (B19), matching the saved view state
$ws.Range("B19").Select()
